# Update price list date and prices on "Hoja1"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the list date (A1) by one month (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 20-31
$ws.Range("D20").Value = 565.303
$ws.Range("D21").Value = 626.903
$ws.Range("D22").Value = 808.479
$ws.Range("D23").Value = 1212.716
$ws.Range("D24").Value = 2425.391
$ws.Range("D25").Value = 2911.734
$ws.Range("D26").Value = 570
$ws.Range("D27").Value = 682
$ws.Range("D28").Value = 569
$ws.Range("D29").Value = 680
$ws.Range("D30").Value = 591
$ws.Range("D31").Value = 700
